# Add an "Image" column (L) to the book-import template header row, and
# bring it into view the way the original author's Excel session ended up
# (selection spanning the full header row through the new column).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell - inherits the header row's style (bold/centered "s=1",
# same as A1) automatically because row 1 already carries that row-default
# format.
$ws.Range("L1").Value = "Image"

# Give the new column a sensible width, matching the other header columns.
$ws.Columns("L").ColumnWidth = 18.3

# Match the author's final on-screen selection/scroll position.
$ws.Range("A1:L1").Select()
